$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rows 112 and 113 (columns F..V) were swapped — the match info (teams,
#    scores, odds, timestamps, url) for "Trabzonspor vs Konyaspor" and
#    "Ankaragucu vs Antalyaspor" change places while columns A:E (index,
#    pais, torneio, temporada, data_partida) stay put.
# ---------------------------------------------------------------------------

$cols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")

$row112 = @()
$row113 = @()
foreach ($col in $cols) {
    $row112 += , $ws.Range($col + "112").Value()
    $row113 += , $ws.Range($col + "113").Value()
}

for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "112").Value = $row113[$i]
    $ws.Range($cols[$i] + "113").Value = $row112[$i]
}

# ---------------------------------------------------------------------------
# 2) Append a new row 153: Gaziantep 2 x 2 Adana Demirspor
# ---------------------------------------------------------------------------

# Pull formatting (styles) from row 152 first, so the new cells land with the
# same cell styles already used throughout the sheet (bold/bordered index in
# column A, date-time format in column E) instead of minting new ones.
$ws.Range("A152").Copy()
$ws.Range("A153").PasteSpecial(-4122)
$ws.Range("E152").Copy()
$ws.Range("E153").PasteSpecial(-4122)

$ws.Range("A153").Value = 152
$ws.Range("B153").Value = "turkey"
$ws.Range("C153").Value = "super-lig"
$ws.Range("D153").Value = "2023-2024"
$ws.Range("E153").Value = 45279.75
$ws.Range("F153").Value = "Gaziantep"
$ws.Range("G153").Value = 2
$ws.Range("H153").Value = "Adana Demirspor"
$ws.Range("I153").Value = 2
$ws.Range("J153").Value = 2.53
$ws.Range("K153").Value = "14/12/2023 09:42"
$ws.Range("L153").Value = 3
$ws.Range("M153").Value = "19/12/2023 17:59"
$ws.Range("N153").Value = 3.67
$ws.Range("O153").Value = "14/12/2023 09:42"
$ws.Range("P153").Value = 3.78
$ws.Range("Q153").Value = "19/12/2023 17:56"
$ws.Range("R153").Value = 2.7
$ws.Range("S153").Value = "14/12/2023 09:42"
$ws.Range("T153").Value = 2.33
$ws.Range("U153").Value = "19/12/2023 17:59"
$ws.Range("V153").Value = "https://www.betexplorer.com/football/turkey/super-lig/gaziantep-adanademirspor/hdMCUdQD/"
